# Append one new Q&A pair (question + sparql_query) to the "Foglio1" sheet,
# right after the current last row (90) -> new row 91.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$question = 'Show me an example of a lemma starting with the prefix "ad"'
$query = "select ?lemma ?label ?prefix where {`n  ?lemma a lila:Lemma;`n    	lila:hasPrefix prefix:5;`n  		rdfs:label ?label.`n  prefix:5 rdfs:label ?prefix`n}limit 1"

$row = 91

# Write column B (sparql_query) before column A (question) so the shared-string
# table gets the query at index 177 and the question at index 178, matching
# the order they were added to the workbook.
$ws.Cells.Item($row, 2).Value = $query
$ws.Cells.Item($row, 1).Value = $question

# Match the formatting used by every other data row: column A wraps with
# top vertical alignment, column B just wraps.
$ws.Cells.Item($row, 1).WrapText = $true
$ws.Cells.Item($row, 1).VerticalAlignment = -4160   # xlTop
$ws.Cells.Item($row, 2).WrapText = $true

# The long multi-line query needs a taller row, same as similar rows above.
$ws.Rows.Item($row).RowHeight = 136

# Leave the view scrolled down near the new row, with the selection where
# the author's cursor ended up after entering the new data.
$excel.ActiveWindow.ScrollRow = 89
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B95").Select()
